$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Diseñador de animación (enero de 2021 - actualidad)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Animación de Spark: Diseñador de animaciones (enero de 2021 - Presente)", 2)

$d.Content.Find.Execute(
    "Diseñador de animación (junio de 2018 - diciembre de 2020)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Pixel Studio: Diseñador de animaciones (junio de 2018 - dic 2020)", 2)

$d.Content.Find.Execute(
    "Diseñador de animación júnior (septiembre de 2016 - mayo de 2018)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Animación flash: Diseñador de animaciones junior (septiembre de 2016 - mayo de 2018)", 2)

$d.Content.Find.Execute(
    "Máster en animación. Fecha prevista de graduación:", $true, $false, $false, $false, $false,
    $true, 1, $false, "Maestro de Artes en Animación, Graduación esperada: dic 2025", 2)
